$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "master table" values in C17:D20 (fSorb (data) / sfSorb (data))
$ws.Range("C17").Value = 0.94362105126578921
$ws.Range("D17").Value = 0.0059433182421027844

$ws.Range("C18").Value = 0.93519060357189066
$ws.Range("D18").Value = 0.010052073987640393

$ws.Range("C19").Value = 0.90146617923715444
$ws.Range("D19").Value = 0.015231309369218682

$ws.Range("C20").Value = 0.9003949709257405
$ws.Range("D20").Value = 0.0058232625164407106

# Move the selection to G15 as recorded in the saved workbook view
[void]$ws.Range("G15").Select()
